$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Vtn"
$ws.Range("C2").Value = "Tnfrsf11b"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.021590333333334
$ws.Range("H2").Value = 6.064771
$ws.Range("I2").Value = 0.01116262347650641
$ws.Range("J2").Value = 0.01116262347650641
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.1176943333333333
$ws.Range("N2").Value = 0.353083
$ws.Range("O2").Value = 0.05829606481791055
$ws.Range("P2").Value = 0.05829606481791055
$ws.Range("Q2").Value = 0.2379297265547778
$ws.Range("R2").Value = 2.141367538993
$ws.Range("S2").Value = 0.0006507370217243478
$ws.Range("T2").Value = 0.0006507370217243479

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Vtn"
$ws.Range("C3").Value = "Tnfrsf11b"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.021590333333334
$ws.Range("H3").Value = 6.064771
$ws.Range("I3").Value = 0.01116262347650641
$ws.Range("J3").Value = 0.01116262347650641
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.183046666666667
$ws.Range("N3").Value = 3.54914
$ws.Range("O3").Value = 0.5859837360842607
$ws.Range("P3").Value = 0.5859837360842608
$ws.Range("Q3").Value = 2.391635705215556
$ws.Range("R3").Value = 21.52472134694
$ws.Range("S3").Value = 0.006541115809265106
$ws.Range("T3").Value = 0.006541115809265107

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Vtn"
$ws.Range("C4").Value = "Tnfrsf11b"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.021590333333334
$ws.Range("H4").Value = 6.064771
$ws.Range("I4").Value = 0.01116262347650641
$ws.Range("J4").Value = 0.01116262347650641
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.718166
$ws.Range("N4").Value = 2.154498
$ws.Range("O4").Value = 0.3557201990978286
$ws.Range("P4").Value = 0.3557201990978286
$ws.Range("Q4").Value = 1.451837443328667
$ws.Range("R4").Value = 13.066536989958
$ws.Range("S4").Value = 0.003970770645516957
$ws.Range("T4").Value = 0.003970770645516957

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Vtn"
$ws.Range("C5").Value = "Tnfrsf11b"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 25.140634
$ws.Range("H5").Value = 75.421902
$ws.Range("I5").Value = 0.1388191398995883
$ws.Range("J5").Value = 0.1388191398995883
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.1176943333333333
$ws.Range("N5").Value = 0.353083
$ws.Range("O5").Value = 0.05829606481791055
$ws.Range("P5").Value = 0.05829606481791055
$ws.Range("Q5").Value = 2.958910158207333
$ws.Range("R5").Value = 26.630191423866
$ws.Range("S5").Value = 0.008092609577552992
$ws.Range("T5").Value = 0.00809260957755299

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Vtn"
$ws.Range("C6").Value = "Tnfrsf11b"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 25.140634
$ws.Range("H6").Value = 75.421902
$ws.Range("I6").Value = 0.1388191398995883
$ws.Range("J6").Value = 0.1388191398995883
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.183046666666667
$ws.Range("N6").Value = 3.54914
$ws.Range("O6").Value = 0.5859837360842607
$ws.Range("P6").Value = 0.5859837360842608
$ws.Range("Q6").Value = 29.74254325158667
$ws.Range("R6").Value = 267.68288926428
$ws.Range("S6").Value = 0.08134575823836442
$ws.Range("T6").Value = 0.08134575823836442

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Vtn"
$ws.Range("C7").Value = "Tnfrsf11b"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 25.140634
$ws.Range("H7").Value = 75.421902
$ws.Range("I7").Value = 0.1388191398995883
$ws.Range("J7").Value = 0.1388191398995883
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.718166
$ws.Range("N7").Value = 2.154498
$ws.Range("O7").Value = 0.3557201990978286
$ws.Range("P7").Value = 0.3557201990978286
$ws.Range("Q7").Value = 18.055148557244
$ws.Range("R7").Value = 162.496337015196
$ws.Range("S7").Value = 0.04938077208367088
$ws.Range("T7").Value = 0.04938077208367087

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Vtn"
$ws.Range("C8").Value = "Tnfrsf11b"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 153.9412893333333
$ws.Range("H8").Value = 461.823868
$ws.Range("I8").Value = 0.8500182366239053
$ws.Range("J8").Value = 0.8500182366239052
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.1176943333333333
$ws.Range("N8").Value = 0.353083
$ws.Range("O8").Value = 0.05829606481791055
$ws.Range("P8").Value = 0.05829606481791055
$ws.Range("Q8").Value = 18.11801742056045
$ws.Range("R8").Value = 163.062156785044
$ws.Range("S8").Value = 0.0495527182186332
$ws.Range("T8").Value = 0.0495527182186332

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Vtn"
$ws.Range("C9").Value = "Tnfrsf11b"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 153.9412893333333
$ws.Range("H9").Value = 461.823868
$ws.Range("I9").Value = 0.8500182366239053
$ws.Range("J9").Value = 0.8500182366239052
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.183046666666667
$ws.Range("N9").Value = 3.54914
$ws.Range("O9").Value = 0.5859837360842607
$ws.Range("P9").Value = 0.5859837360842608
$ws.Range("Q9").Value = 182.1197292081689
$ws.Range("R9").Value = 1639.07756287352
$ws.Range("S9").Value = 0.4980968620366312
$ws.Range("T9").Value = 0.4980968620366312

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Vtn"
$ws.Range("C10").Value = "Tnfrsf11b"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 153.9412893333333
$ws.Range("H10").Value = 461.823868
$ws.Range("I10").Value = 0.8500182366239053
$ws.Range("J10").Value = 0.8500182366239052
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.718166
$ws.Range("N10").Value = 2.154498
$ws.Range("O10").Value = 0.3557201990978286
$ws.Range("P10").Value = 0.3557201990978286
$ws.Range("Q10").Value = 110.5553999953627
$ws.Range("R10").Value = 994.9985999582639
$ws.Range("S10").Value = 0.3023686563686408
$ws.Range("T10").Value = 0.3023686563686408

Write-Host "Edit complete"
